$d = $word.ActiveDocument

function Insert-RunsXml([object]$range, [string]$innerXml) {
    # Wrap a fragment of <w:r>/<w:bookmarkStart>/... elements in a minimal
    # WordProcessingML package so Range.InsertXML can splice it into the
    # run stream at $range (after $range's own content has been removed).
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1. "Implemente um programa..." paragraph: re-typing the same text over
#    the span that held the _GoBack bookmark merges the two runs back
#    into one and drops the now-stale bookmark, with no net change in
#    the visible text.
# ---------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute(
    "Implemente um programa que dado 02 valores atribuídos as variáveis A e B consecutivamente, efetuar a troca dos valores das mesmas de forma que A passe a ter o valor de B e B o valor de A.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implemente um programa que dado 02 valores atribuídos as variáveis A e B consecutivamente, efetuar a troca dos valores das mesmas de forma que A passe a ter o valor de B e B o valor de A.",
    2
)

# ---------------------------------------------------------------------
# 2. "Desenvolva um programa..." paragraph gets substantially reworded,
#    in several separate edit passes (hence the many runs in the final
#    XML, including the _GoBack bookmark moving to a new spot).
# ---------------------------------------------------------------------
$r2 = $d.Content
$ok2 = $r2.Find.Execute(
    "Desenvolva um programa para que dado um intervalo de números inteiros, seu programa apresente os quadrados dos números do intervalo. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$r2.Delete()
$r2ins = $d.Range($r2.Start, $r2.Start)

$p2xml = (
    '<w:r><w:t>Desenvolva um programa para que dado um intervalo de números inteiros</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>S</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">eu programa </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">deve </w:t></w:r>' +
    '<w:r><w:t>a</w:t></w:r>' +
    '<w:r><w:t>present</w:t></w:r>' +
    '<w:r><w:t>ar</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">a potenciação ao quadrado dos números </w:t></w:r>' +
    '<w:r><w:t>do int</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>ervalo</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> informado</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>'
)
Insert-RunsXml $r2ins $p2xml

# ---------------------------------------------------------------------
# 3. "Desenvolver uma classe..." -> "Desenvolver um programa..."
# ---------------------------------------------------------------------
$r3 = $d.Content
$ok3 = $r3.Find.Execute(
    "Desenvolver uma classe que imprima os números perfeitos compreendidos entre 1 e 500",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
$r3.Delete()
$r3ins = $d.Range($r3.Start, $r3.Start)

$p3xml = (
    '<w:r><w:t xml:space="preserve">Desenvolver um programa </w:t></w:r>' +
    '<w:r><w:t>que imprima os números perfeitos compreendidos entre 1 e 500</w:t></w:r>'
)
Insert-RunsXml $r3ins $p3xml
